# Train Model for operation 1
# Falta el modulo de pronostico para el modo de operacion 1
#
# Adds a new "Ea_train" variable row to the InfluxDBVariables sheet,
# right after the other "DTPlantaBiogas" device rows and before the
# "Modulo solar-eolico" block (i.e. as a new row 158, pushing the
# existing rows 158-206 down to 159-207).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InfluxDBVariables")
$ws.Activate()

# Insert a new blank row at row 158 (shifts rows 158:206 -> 159:207).
$ws.Rows.Item(158).Insert()

# Populate the new row: Device / Name / Tag.
$ws.Cells.Item(158, 1).Value = "DTPlantaBiogas"
$ws.Cells.Item(158, 2).Value = "Ea_train"
$ws.Cells.Item(158, 3).Value = "Ea_train"

# Match the author's final view state (scrolled/selected cell).
$ws.Range("C158").Select()
